$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.291"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05786"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.484"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.331"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8088"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8764"
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1384"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07278"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03091"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03057"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09314"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.861"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001542"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04709"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006029"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006022"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001291"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'0.00008699"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'2.143"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3181"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Value = "'0.0002349"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.03772"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1053"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("D42").Value = "'0.002540"
$ws.Range("D42").Style = "Normal"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003177"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").Value = "'0.007640"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "'0.00005478"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.5898"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002528"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
